$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.188.70'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +8.17%  '
$ws.Range('E2').ClearFormats()

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.877.60'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.52%  '
$ws.Range('E3').ClearFormats()

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E4').ClearFormats()

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.49'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.20%  '
$ws.Range('E5').ClearFormats()

$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9992'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E6').ClearFormats()

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4989'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.42%  '
$ws.Range('E7').ClearFormats()

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '45.74'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +9.20%  '
$ws.Range('E8').ClearFormats()

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2858'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.57%  '
$ws.Range('E9').ClearFormats()

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06614'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.08%  '
$ws.Range('E10').ClearFormats()

$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.870.04'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.08%  '
$ws.Range('E11').ClearFormats()

$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '17.11'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +5.01%  '
$ws.Range('E12').ClearFormats()

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07215'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('E13').ClearFormats()

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6641'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +7.97%  '
$ws.Range('E14').ClearFormats()

$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '85.69'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +7.95%  '
$ws.Range('E15').ClearFormats()

$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.833'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.61%  '
$ws.Range('E16').ClearFormats()

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.172.24'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +8.21%  '
$ws.Range('E17').ClearFormats()

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9990'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('E18').ClearFormats()

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.01'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +10.76%  '
$ws.Range('E19').ClearFormats()

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007578'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +5.38%  '
$ws.Range('E20').ClearFormats()

$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9981'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('E21').ClearFormats()

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.108.53'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.94%  '
$ws.Range('E22').ClearFormats()

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.789'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.64%  '
$ws.Range('E23').ClearFormats()

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.086'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.23%  '
$ws.Range('E24').ClearFormats()

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.527'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +6.80%  '
$ws.Range('E25').ClearFormats()

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '144.88'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.08%  '
$ws.Range('E26').ClearFormats()

$ws.Range('B27').Value = 'BitcoinCash'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '136.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +24.96%  '
$ws.Range('E27').ClearFormats()

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.82'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +8.19%  '
$ws.Range('E28').ClearFormats()

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.960'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.80%  '
$ws.Range('E29').ClearFormats()

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.392'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E30').ClearFormats()

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.262'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +5.29%  '
$ws.Range('E31').ClearFormats()

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08645'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.77%  '
$ws.Range('E32').ClearFormats()

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.915'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.08%  '
$ws.Range('E33').ClearFormats()

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05096'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +8.40%  '
$ws.Range('E34').ClearFormats()

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.135'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.02%  '
$ws.Range('E35').ClearFormats()

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6904'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +9.18%  '
$ws.Range('E36').ClearFormats()

$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.696'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.52%  '
$ws.Range('E37').ClearFormats()

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.361'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +15.58%  '
$ws.Range('E38').ClearFormats()

$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.748'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.15%  '
$ws.Range('E39').ClearFormats()

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9627'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.83%  '
$ws.Range('E40').ClearFormats()

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01636'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +7.01%  '
$ws.Range('E41').ClearFormats()

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.135'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.15%  '
$ws.Range('E42').ClearFormats()

$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '103.92'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.90%  '
$ws.Range('E43').ClearFormats()

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9992'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E44').ClearFormats()

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4213'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.60%  '
$ws.Range('E45').ClearFormats()

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.524'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.43%  '
$ws.Range('E46').ClearFormats()

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1262'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.03%  '
$ws.Range('E47').ClearFormats()

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05641'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.45%  '
$ws.Range('E48').ClearFormats()

$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.71'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +7.98%  '
$ws.Range('E49').ClearFormats()

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.300'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.15%  '
$ws.Range('E50').ClearFormats()

$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3747'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +8.84%  '
$ws.Range('E51').ClearFormats()
